$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "37.219.84"
Set-TextValue 2 5 "  +1.44%  "
Set-TextValue 3 4 "2.030.58"
Set-TextValue 3 5 "  +3.47%  "
Set-TextValue 4 5 "  -0.08%  "
Set-TextValue 5 4 "247.85"
Set-TextValue 5 5 "  +1.32%  "
Set-TextValue 6 5 "  +1.86%  "
Set-TextValue 7 4 "60.59"
Set-TextValue 7 5 "  -1.45%  "
Set-TextValue 8 5 "  -0.02%  "
Set-TextValue 9 4 "0.396"
Set-TextValue 9 5 "  +5.05%  "
Set-TextValue 10 5 "  +2.14%  "
Set-TextValue 11 5 "  +1.93%  "
Set-TextValue 12 4 "15.27"
Set-TextValue 12 5 "  +6.54%  "
Set-TextValue 13 4 "0.864"
Set-TextValue 13 5 "  +3.65%  "
Set-TextValue 14 4 "22.57"
Set-TextValue 14 5 "  +2.32%  "
Set-TextValue 15 4 "2.325.96"
Set-TextValue 15 5 "  +3.31%  "
Set-TextValue 16 4 "5.52"
Set-TextValue 16 5 "  +4.08%  "
Set-TextValue 17 4 "2.031.09"
Set-TextValue 17 5 "  +3.83%  "
Set-TextValue 18 4 "37.155.79"
Set-TextValue 18 5 "  +1.47%  "
Set-TextValue 19 4 "70.76"
Set-TextValue 19 5 "  +1.28%  "
Set-TextValue 20 5 "  +1.44%  "
Set-TextValue 21 4 "5.27"
Set-TextValue 21 5 "  +3.56%  "
Set-TextValue 22 4 "231.44"
Set-TextValue 22 5 "  +0.50%  "
Set-TextValue 23 5 "  +0.02%  "
Set-TextValue 24 5 "  +2.35%  "
Set-TextValue 25 5 "  +0.65%  "
Set-TextValue 26 4 "9.48"
Set-TextValue 26 5 "  +2.64%  "
Set-TextValue 27 4 "164.24"
Set-TextValue 27 5 "  +2.11%  "
Set-TextValue 28 5 "  -2.95%  "
Set-TextValue 29 4 "19.89"
Set-TextValue 29 5 "  +2.39%  "
Set-TextValue 30 5 "  +7.97%  "
Set-TextValue 31 5 "  +2.10%  "
Set-TextValue 32 4 "4.86"
Set-TextValue 32 5 "  +1.41%  "
Set-TextValue 33 4 "0.0676"
Set-TextValue 33 5 "  +9.24%  "
Set-TextValue 34 4 "4.55"
Set-TextValue 34 5 "  +2.30%  "
Set-TextValue 35 4 "2.50"
Set-TextValue 35 5 "  +9.88%  "
Set-TextValue 36 4 "3.46"
Set-TextValue 36 5 "  -3.39%  "
Set-TextValue 37 5 "  -0.20%  "
Set-TextValue 38 5 "  +2.02%  "
Set-TextValue 39 5 "  -1.53%  "
Set-TextValue 40 4 "0.0983"
Set-TextValue 40 5 "  +0.07%  "
Set-TextValue 41 5 "  +1.45%  "
Set-TextValue 42 5 "  +2.03%  "
Set-TextValue 43 2 "InjectiveProtocol"
Set-TextValue 43 3 "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue 43 4 "16.97"
Set-TextValue 43 5 "  +5.56%  "
Set-TextValue 44 2 "VeChain"
Set-TextValue 44 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 44 4 "0.0215"
Set-TextValue 44 5 "  +1.93%  "
Set-TextValue 45 4 "92.52"
Set-TextValue 45 5 "  +4.01%  "
Set-TextValue 46 2 "ARBITRUM"
Set-TextValue 46 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue 46 4 "1.07"
Set-TextValue 46 5 "  +3.57%  "
Set-TextValue 47 2 "Maker"
Set-TextValue 47 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue 47 4 "1.391.23"
Set-TextValue 47 5 "  +1.40%  "
Set-TextValue 48 4 "7.56"
Set-TextValue 48 5 "  +5.40%  "
Set-TextValue 49 4 "2.15"
Set-TextValue 49 5 "  +16.78%  "
Set-TextValue 50 4 "2.86"
Set-TextValue 50 5 "  +0.61%  "
Set-TextValue 51 4 "46.80"
Set-TextValue 51 5 "  +3.20%  "

Write-Host "Applied cryptos list update"
